$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '64.465.73'
$ws.Range("D3").Value = '3.358.20'
$ws.Range("E3").Value = '  -2.23%  '
$ws.Range("E4").Value = '  +0.00%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '555.93'
$ws.Range("E5").Value = '  -2.72%  '
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '175.40'
$ws.Range("E6").Value = '  +0.53%  '
$ws.Range("E7").Value = '  -0.55%  '
$ws.Range("D8").Value = '3.348.92'
$ws.Range("E8").Value = '  -2.31%  '
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '1.00'
$ws.Range("E9").Value = '  +0.02%  '
$ws.Range("E10").Value = '  +3.45%  '
$ws.Range("E11").Value = '  +0.92%  '
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '54.48'
$ws.Range("E12").Value = '  -0.92%  '
$ws.Range("E13").Value = '  +1.07%  '
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '9.09'
$ws.Range("E14").Value = '  -0.24%  '
$ws.Range("D15").Value = '3.887.99'
$ws.Range("E15").Value = '  -2.39%  '
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '18.42'
$ws.Range("E16").Value = '  +1.98%  '
$ws.Range("E17").Value = '  -1.79%  '
$ws.Range("D18").Value = '3.356.53'
$ws.Range("E18").Value = '  -2.49%  '
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '11.85'
$ws.Range("E19").Value = '  +0.09%  '
$ws.Range("D20").Value = '64.396.69'
$ws.Range("E20").Value = '  -0.65%  '
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '0.985'
$ws.Range("E21").Value = '  -0.24%  '
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '461.86'
$ws.Range("E22").Value = '  +13.43%  '
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '4.85'
$ws.Range("E23").Value = '  +11.54%  '
$ws.Range("E24").Value = '  -2.25%  '
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '86.19'
$ws.Range("E25").Value = '  +3.23%  '
$ws.Range("E26").Value = '  +1.07%  '
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '10.96'
$ws.Range("E27").Value = '  +1.42%  '
$ws.Range("E28").Value = '  +2.00%  '
$ws.Range("E29").Value = '  -1.72%  '
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '30.10'
$ws.Range("E30").Value = '  +0.95%  '
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '6.67'
$ws.Range("E31").Value = '  +1.16%  '
$ws.Range("B32").Value = 'Cosmos'
$ws.Range("C32").Value = 'https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom'
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '11.47'
$ws.Range("E32").Value = '  -0.31%  '
$ws.Range("B33").Value = 'Bittensor'
$ws.Range("C33").Value = 'https://coinranking.com/coin/pgv7xSFi6+bittensor-tao'
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '583.00'
$ws.Range("E33").Value = '  -0.57%  '
$ws.Range("E34").Value = '  -0.13%  '
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '58.81'
$ws.Range("E35").Value = '  -1.35%  '
$ws.Range("E36").Value = '  +0.07%  '
$ws.Range("E37").Value = '  -8.24%  '
$ws.Range("E38").Value = '  -0.08%  '
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '35.67'
$ws.Range("E39").Value = '  -1.28%  '
$ws.Range("D40").Value = '0.0₃0758'
$ws.Range("E40").Value = '  -1.24%  '
$ws.Range("E41").Value = '  -0.18%  '
$ws.Range("D42").Value = '3.096.56'
$ws.Range("E42").Value = '  -2.69%  '
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '0.999'
$ws.Range("E43").Value = '  -0.06%  '
$ws.Range("B44").Value = 'ThetaToken'
$ws.Range("C44").Value = 'https://coinranking.com/coin/B42IRxNtoYmwK+thetatoken-theta'
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '2.80'
$ws.Range("E44").Value = '  -4.11%  '
$ws.Range("B45").Value = 'Fetch.AI'
$ws.Range("C45").Value = 'https://coinranking.com/coin/AWma-WzFHmKVQ+fetchai-fet'
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '2.52'
$ws.Range("E45").Value = '  +0.89%  '
$ws.Range("E46").Value = '  +0.76%  '
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '3.20'
$ws.Range("E47").Value = '  -1.74%  '
$ws.Range("E48").Value = '  +0.38%  '
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '2.58'
$ws.Range("E49").Value = '  -1.73%  '
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '8.39'
$ws.Range("E50").Value = '  -0.53%  '
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '135.38'
$ws.Range("E51").Value = '  -1.89%  '
